# tracker.xlsx - add "Product of Array Except Itself" leetcode entry (row 32)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fill in the row's content first (order matters for the shared
#        string table: url text before notes text). ---
$ws.Range("A32").Value = 238
$ws.Range("B32").Value = "Medium"
$ws.Hyperlinks.Add($ws.Range("D32"), "https://leetcode.com/problems/product-of-array-except-self/description/")
$ws.Range("F32").Value = "O(n)"
$ws.Range("G32").Value = "Initialise the res array, left pass on it using premul variable and then overwrite over it with right pass using postmul."
$ws.Range("H32").Value = 45503
$ws.Range("K32").Value = "Medium"

# --- 2. Copy the formatting (fill/border/number format/hyperlink font)
#        of the already fully-populated row 28 onto row 32 last, so every
#        style lands on the same shared xf row 28 uses (incl. undoing the
#        one-off hyperlink style Excel mints when a link is added). ---
$ws.Range("A28:K28").Copy()
$ws.Range("A32:K32").PasteSpecial(-4122)   # xlPasteFormats

# --- 3. Match the sheet's scroll position / active selection. ---
$ws.Range("C30").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
